# Update "Fruta, Vega Modelo de Temuco - Guayaba" weekly price sheet.
# The rows for this market/product got re-paired with (different) dates,
# volumes and prices -- effectively a re-shuffle of rows 3-13 across the
# Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken from the target workbook state.
$rows = @(
    @{ Row = 3;  D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 4;  D = 44418; M = 40;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 5;  D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 6;  D = 44476; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 7;  D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 8;  D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 9;  D = 44432; M = 30;  N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 10; D = 44473; M = 120; N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 11; D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 12; D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 13; D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 19).Value = $r.S   # S: Precio $/Kg
}
